# ---------------------------------------------------------------------------
# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. before the
#    worksheet that is currently named "2021-Q2"), and populate it with the
#    new fund-holding data.
# 2. Update the "总计" (summary) sheet: the existing rows shift down by one
#    and a new row for "2022-Q4" is inserted at the top of the data.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: add the "2022-Q4" worksheet by duplicating the existing "2021-Q2"
# sheet (so sheet-level properties such as outline settings come along for
# free) and dropping the duplicate in right before it, then overwrite its
# cell contents with the new quarter's numbers.
# ---------------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2021-Q2")
$srcSheet.Copy($srcSheet) | Out-Null
$q4sheet = $wb.Worksheets.Item("2021-Q2 (2)")
$q4sheet.Name = "2022-Q4"

# A "plain" (unstyled) cell used purely as a format donor below, to strip the
# quote-prefix marker Excel leaves behind after typing a leading "'" into a
# cell (see helper below) -- C2 on the total sheet has no explicit style.
$plainCell = $wb.Worksheets.Item("总计").Range("C2")

function Set-TextValue($range, [string]$text, $donor) {
    # Force $text (which may look like a number) to be stored as a literal
    # text value rather than being parsed into a numeric cell, the way
    # typing a leading apostrophe in Excel does -- then strip the resulting
    # "quote prefix" cell style back off so the cell's format matches a
    # normal, never-touched text cell.
    $range.Value = "'" + $text
    $donor.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# Header row (only D1 differs from the source sheet: 基金金额 -> 基金规模)
$q4sheet.Range("D1").Value = "基金规模"

# Row 2
$q4sheet.Range("A2").Value = 0
Set-TextValue $q4sheet.Range("B2") "006890" $plainCell
$q4sheet.Range("C2").Value = "上投摩根领先优选混合A"
Set-TextValue $q4sheet.Range("D2") "0.29" $plainCell
Set-TextValue $q4sheet.Range("E2") "81.41" $plainCell
Set-TextValue $q4sheet.Range("F2") "3.62" $plainCell
Set-TextValue $q4sheet.Range("G2") "0.0105" $plainCell
$q4sheet.Range("H2").Value = 2

# Row 3
$q4sheet.Range("A3").Value = 1
Set-TextValue $q4sheet.Range("B3") "017098" $plainCell
$q4sheet.Range("C3").Value = "上投摩根领先优选混合C"
Set-TextValue $q4sheet.Range("D3") "0.00" $plainCell
Set-TextValue $q4sheet.Range("E3") "81.41" $plainCell
Set-TextValue $q4sheet.Range("F3") "3.62" $plainCell
$q4sheet.Range("G3").Value = 0
$q4sheet.Range("H3").Value = 2

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 2: update the "总计" summary sheet -- shift the three existing data
# rows down by one (carrying their formatting with them, row by row, from
# bottom to top so nothing gets overwritten before it is copied), then write
# in the new row of values (including the brand-new "2022-Q4" row at the
# top).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A4:D4").Copy() | Out-Null
$totalSheet.Range("A5:D5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$totalSheet.Range("A3:D3").Copy() | Out-Null
$totalSheet.Range("A4:D4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$totalSheet.Range("A2:D2").Copy() | Out-Null
$totalSheet.Range("A3:D3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Row 5 (was row 4): 2020-Q4
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2020-Q4"
$totalSheet.Range("C5").Value = 4
$totalSheet.Range("D5").Value = 0.59

# Row 4 (was row 3): 2021-Q1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q1"
$totalSheet.Range("C4").Value = 12
$totalSheet.Range("D4").Value = 6.33

# Row 3 (was row 2): 2021-Q2
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q2"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.33

# Row 2 (new): 2022-Q4
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.01
